# Update Oracle jobs data
# The job listing rows (2-7) shift down by one: each row's Title/Date/Apply link
# becomes the previous row's old values, and a brand-new posting (job/31390,
# "Agri Finance and R&D Specialist (Open to external applicants)") is bumped to
# the top with a new posting date of 02/05/2026.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$jobs = @(
    @{ Row = 2; Title = "Agri Finance and R&D Specialist (Open to external applicants)"; Date = "02/05/2026"; JobId = "31390" },
    @{ Row = 3; Title = "Specialist – Operations, Fund Risk Management and Oversight (Open to all applicants)"; Date = "02/04/2026"; JobId = "28371" },
    @{ Row = 4; Title = "Human Resources Associate"; Date = "02/01/2026"; JobId = "31613" },
    @{ Row = 5; Title = "Human Resources Associate"; Date = "02/01/2026"; JobId = "31614" },
    @{ Row = 6; Title = "Finance Analyst"; Date = "01/26/2026"; JobId = "31468" },
    @{ Row = 7; Title = "Project Accounting & Financial Management Officer"; Date = "01/26/2026"; JobId = "31469" }
)

foreach ($job in $jobs) {
    $r = $job.Row

    $ws.Range("B$r").Value = $job.Title

    # Force the posting date to be stored as literal text (matching the
    # original inline-string cells) instead of being auto-parsed into a
    # date serial number. Applying a temporary text number format prevents
    # Excel's date auto-detection, then the style is reset to Normal so no
    # stray cell formatting is left behind.
    $ws.Range("D$r").NumberFormat = "@"
    $ws.Range("D$r").Value = $job.Date
    $ws.Range("D$r").Style = "Normal"

    $url = "https://estm.fa.em2.oraclecloud.com/hcmUI/CandidateExperience/en/sites/CX_1/job/" + $job.JobId + "/?location=India&locationId=300000000440677&locationLevel=country&mode=location"
    $ws.Range("E$r").Formula = '=HYPERLINK("' + $url + '", "Apply")'
}
